$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.625.69"
$ws.Range("E2").Value = "  +4.05%  "
$ws.Range("D3").Value = "2.468.25"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'323.33"
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("D6").Value = "'105.18"
$ws.Range("E6").Value = "  +3.23%  "
$ws.Range("E7").Value = "  +1.41%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").Value = "'36.08"
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("D11").Value = "'0.0815"
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "'18.29"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").Value = "'7.08"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").Value = "2.853.26"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").Value = "2.469.75"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "'0.843"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "46.495.42"
$ws.Range("E18").Value = "  +4.17%  "
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Value = "'6.46"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").Value = "0.0₃0936"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").Value = "'70.47"
$ws.Range("E22").Value = "  +2.40%  "
$ws.Range("D23").Value = "'249.01"
$ws.Range("E23").Value = "  +2.75%  "
$ws.Range("E24").Value = "  +4.40%  "
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("E26").Value = "  +3.35%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'2.30"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("E29").Value = "  +3.19%  "
$ws.Range("D30").Value = "'35.20"
$ws.Range("E30").Value = "  +4.39%  "
$ws.Range("D31").Value = "'49.62"
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("E32").Value = "  +2.14%  "
$ws.Range("D33").Value = "'19.62"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("D37").Value = "'4.63"
$ws.Range("E37").Value = "  +3.25%  "
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").Value = "'2.94"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("D40").Value = "'123.22"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").Value = "'2.23"
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("D43").Value = "'20.64"
$ws.Range("E43").Value = "  -2.05%  "
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D45").Value = "1.981.82"
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("E48").Value = "  +3.77%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'5.35"
$ws.Range("E49").Value = "  +15.98%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'8.94"
$ws.Range("E50").Value = "  -4.60%  "
$ws.Range("D51").Value = "'79.36"
$ws.Range("E51").Value = "  +4.97%  "
